# edit.ps1 - apply the "updated readme and ppt" commit to TranBrian-Presentation.pptx
#
# Summary of changes:
#  1. Slide 1 (title slide): split the "Web Dev w/ .NET Final Project" run into two
#     runs, and add a new paragraph "Professor Konstantopou" to the subtitle.
#  2. Slide 8 (Demo!): merge the "Collision Detection" + "?" runs into a single run,
#     and split "Jumping Data" into a curly-quoted "Jumping" run + "Data" run.
#  3. Add a new slide 9 ("References") with a content placeholder containing three
#     hyperlinked bullet points.

$p = $ppt.ActivePresentation

# -------------------------------------------------------------------------
# 1. Slide 1 - title slide subtitle text
# -------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2).TextFrame.TextRange

# Split "Web Dev w/ .NET Final Project" (paragraph 3) into two runs:
#   "Web Dev w/ .NET Final " + "Project"
$para3 = $subtitle.Paragraphs(3)
$splitPoint = "Web Dev w/ .NET Final ".Length
$tail = $para3.Characters($splitPoint + 1, $para3.Text.Length - $splitPoint)
$tail.Text = "Project"

# Add a new paragraph after it: "Professor " + "Konstantopou"
$subtitle.InsertAfter("`rProfessor Konstantopou")

$para4 = $slide1.Shapes.Item(2).TextFrame.TextRange.Paragraphs(4)
$namePart = $para4.Characters("Professor ".Length + 1, "Konstantopou".Length)
$namePart.Text = "Konstantopou"

# -------------------------------------------------------------------------
# 2. Slide 8 - Demo! bullets
# -------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(2).TextFrame.TextRange

# Merge "Collision Detection" + "?" into a single run "Collision Detection?"
$body8.Paragraphs(6).Text = "Collision Detection?"

# Split "Jumping Data" into a curly-quoted "Jumping" run + "Data" run
$jumpPara = $body8.Paragraphs(9)
$jumpPara.Text = [char]0x201C + "Jumping" + [char]0x201D + " Data"
$dataRun = $jumpPara.Characters(11, 4)
$dataRun.Text = "Data"

# -------------------------------------------------------------------------
# 3. New slide 9 - References
# -------------------------------------------------------------------------
$slide9 = $p.Slides.Add($p.Slides.Count + 1, 2)

$slide9.Shapes.Item(1).TextFrame.TextRange.Text = "References"

$refs = $slide9.Shapes.Item(2).TextFrame.TextRange
$refs.Text = "https://www.asp.net/signalr`rhttps://github.com/NTaylorMullen/ShootR/tree/master/ShootR`rhttps://phaser.io/docs/2.6.2/index"

$refs.Paragraphs(1).ActionSettings(1).Hyperlink.Address = "https://www.asp.net/signalr"
$refs.Paragraphs(2).ActionSettings(1).Hyperlink.Address = "https://www.asp.net/signalr"
$refs.Paragraphs(3).ActionSettings(1).Hyperlink.Address = "https://phaser.io/docs/2.6.2/index"

Write-Output "done"
